$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16 (pushes the "Relay States" / clustering /
# time-series block down by one row) to make room for the new
# "Support Vector Classifier (SVC)" / "Binary Classification" /
# "Industrial System Monitoring" entry.
$ws.Rows.Item(16).Insert()

# Populate the new row's values. Set C16 before A16 so the shared-string
# table gains "Industrial System Monitoring" before "Support Vector
# Classifier (SVC)" (matches the authoring order in the target file).
$ws.Range("C16").Value() = "Industrial System Monitoring"
$ws.Range("B16").Value() = "Binary Classification"
$ws.Range("A16").Value() = "Support Vector Classifier (SVC)"

# Give the new row a hyperlink on the project-name cell, consistent with
# the other rows in the table.
$ws.Hyperlinks.Add($ws.Range("C16"), "c. Jupyter Notebooks\Industrial System Monitoring.ipynb")

# Restore the "plain" (non auto-hyperlink-styled) look used throughout
# this table by re-applying the formatting from an existing row of the
# same visual pattern (unmerged A/B/C, single-line row).
$ws.Range("A9:C9").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the saved selection/scroll position recorded for this edit.
$ws.Range("A17").Select()
$ws.Application.ActiveWindow.ScrollRow = 8
